$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column S immediately after the last data column (R).
# Inserting here takes on the formatting of the column to its left (R),
# which reproduces the exact per-row cell styles (s="...") the target
# workbook uses for the new 2022 column.
$ws.Columns("S").Insert()

# Header (row 4 holds the year labels) — new year column.
$ws.Range("S4").Value = 2022

# Data rows 5-34: figures for 2022, matching the existing 2007-2021 columns.
$ws.Range("S5").Value = 135
$ws.Range("S6").Value = 99
$ws.Range("S7").Value = 36
$ws.Range("S8").Value = 97
$ws.Range("S9").Value = 80
$ws.Range("S10").Value = 17
$ws.Range("S11").Value = 17
$ws.Range("S12").Value = 11
$ws.Range("S13").Value = 6
$ws.Range("S14").Value = 5
$ws.Range("S15").Value = 3
$ws.Range("S16").Value = 2
$ws.Range("S17").Value = "-"
$ws.Range("S18").Value = "-"
$ws.Range("S19").Value = "-"
$ws.Range("S20").Value = 6
$ws.Range("S21").Value = 1
$ws.Range("S22").Value = 5
$ws.Range("S23").Value = "-"
$ws.Range("S24").Value = "-"
$ws.Range("S25").Value = "-"
$ws.Range("S26").Value = 10
$ws.Range("S27").Value = 4
$ws.Range("S28").Value = 6
$ws.Range("S29").Value = "-"
$ws.Range("S30").Value = "-"
$ws.Range("S31").Value = "-"
$ws.Range("S32").Value = "-"
$ws.Range("S33").Value = "-"
$ws.Range("S34").Value = "-"

# Match the selection recorded in the saved workbook.
$ws.Range("S3").Select()
